$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.924350738525391
$ws.Range("B1").Value = 4.765216827392578
$ws.Range("C1").Value = 3.463376522064209
$ws.Range("D1").Value = 2.608125686645508
$ws.Range("E1").Value = 1.927480936050415
